# Commit message: "Updated the file old date referances."
#
# The "Screenshots" sheet has a question (cell A7) that referenced a
# specific/stale date range ("September 8th to September 12th"); update it
# to a non-dated phrase ("the first full week of the semester"), and leave
# the selection on that cell (matches the author re-selecting A7 after the
# edit, instead of the old A5).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Screenshots")

$ws.Range("A7").Value = '3. On our Class Website, navigate to the "Office Hours" tab and take a screenshot of the office hours for the first full week of the semester.'

$ws.Range("A7").Select()
